$wb = $excel.ActiveWorkbook

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 18204.445
$ws.Range("J69").Value = 17475.834
$ws.Range("L69").Value = 52427.50199999999
$ws.Range("N69").Value = -54175.50199999999

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 18204.445
$ws.Range("J72").Value = 17475.834
$ws.Range("L72").Value = 157282.506
$ws.Range("N72").Value = -166018.506

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4237.6
$ws.Range("I76").Value = 3779
$ws.Range("J76").Value = 4696.2
$ws.Range("K76").Value = 3779
$ws.Range("L76").Value = 4696.2
$ws.Range("M76").Value = -3464
$ws.Range("N76").Value = -5326.2

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 4237.6
$ws.Range("I79").Value = 3779
$ws.Range("J79").Value = 4696.2
$ws.Range("K79").Value = 3779
$ws.Range("L79").Value = 4696.2
$ws.Range("M79").Value = -2687
$ws.Range("N79").Value = -6880.2

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1454.7727
$ws.Range("I98").Value = 1200.25
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1200.25
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = 297.75
$ws.Range("N98").Value = -6996

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1454.7727
$ws.Range("I122").Value = 1200.25
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3600.75
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1150.75
$ws.Range("N122").Value = -16900

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 8072.6665
$ws.Range("J138").Value = 8660.821
$ws.Range("L138").Value = 25982.463
$ws.Range("N138").Value = -36262.463

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 432.85715
$ws.Range("I5").Value = 432.85715
$ws.Range("K5").Value = 432.85715
$ws.Range("M5").Value = -320.85715

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5837.706
$ws.Range("I61").Value = 5837.706
$ws.Range("K61").Value = 5837.706
$ws.Range("M61").Value = -5625.706

# ARM row 69
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H69").Value = 395000
$ws.Range("I69").Value = 395000
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 395000
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -394251
$ws.Range("N69").ClearContents()

# ARM row 72
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H72").Value = 395000
$ws.Range("I72").Value = 395000
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 1185000
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -1181256
$ws.Range("N72").ClearContents()

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 89086080
$ws.Range("I74").Value = 152715970
$ws.Range("J74").Value = 4237
$ws.Range("K74").Value = 152715970
$ws.Range("L74").Value = 4237
$ws.Range("M74").Value = -152715096
$ws.Range("N74").Value = -5985

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 89086080
$ws.Range("I77").Value = 152715970
$ws.Range("J77").Value = 4237
$ws.Range("K77").Value = 763579850
$ws.Range("L77").Value = 21185
$ws.Range("M77").Value = -763575482
$ws.Range("N77").Value = -29921

# ARM row 97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 762.4375
$ws.Range("J97").Value = 1015.3333
$ws.Range("L97").Value = 1015.3333
$ws.Range("N97").Value = -2007.3333

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5837.706
$ws.Range("I136").Value = 5837.706
$ws.Range("K136").Value = 17513.118
$ws.Range("M136").Value = -14963.118

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 432.85715
$ws.Range("I4").Value = 432.85715
$ws.Range("K4").Value = 432.85715
$ws.Range("M4").Value = -317.85715

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 66706668
$ws.Range("I86").Value = 97502.5
$ws.Range("J86").Value = 111112780
$ws.Range("K86").Value = 97502.5
$ws.Range("L86").Value = 111112780
$ws.Range("M86").Value = -96379.5
$ws.Range("N86").Value = -111115026

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 66706668
$ws.Range("I89").Value = 97502.5
$ws.Range("J89").Value = 111112780
$ws.Range("K89").Value = 487512.5
$ws.Range("L89").Value = 555563900
$ws.Range("M89").Value = -481896.5
$ws.Range("N89").Value = -555575132

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1485.0526
$ws.Range("I105").Value = 1428.2667
$ws.Range("K105").Value = 1428.2667
$ws.Range("M105").Value = 318.7333000000001

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 83333640
$ws.Range("I7").Value = 125000330
$ws.Range("J7").Value = 263
$ws.Range("K7").Value = 125000330
$ws.Range("L7").Value = 263
$ws.Range("M7").Value = -125000217
$ws.Range("N7").Value = -489

# CRP row 22
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 558.9
$ws.Range("I22").Value = 837.25
$ws.Range("J22").Value = 373.33334
$ws.Range("K22").Value = 837.25
$ws.Range("L22").Value = 373.33334
$ws.Range("M22").Value = -487.25
$ws.Range("N22").Value = -1073.33334

# CRP row 57
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 962.5
$ws.Range("I57").Value = 962.5
$ws.Range("K57").Value = 962.5
$ws.Range("M57").Value = -402.5

# CRP row 74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# CRP row 77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# CUL row 11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 1125
$ws.Range("I11").Value = 800
$ws.Range("J11").Value = 1233.3334
$ws.Range("K11").Value = 2400
$ws.Range("L11").Value = 3700.0002
$ws.Range("M11").Value = -2260
$ws.Range("N11").Value = -3980.0002

# CUL row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 895.1
$ws.Range("I14").Value = 895.1
$ws.Range("K14").Value = 2685.3
$ws.Range("M14").Value = -2512.3

# CUL row 33
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 278.57144
$ws.Range("I33").Value = 271.66666
$ws.Range("K33").Value = 1629.99996
$ws.Range("M33").Value = -1346.99996

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3101
$ws.Range("I113").Value = 2999.25
$ws.Range("J113").Value = 3128.1333
$ws.Range("K113").Value = 8997.75
$ws.Range("L113").Value = 9384.3999
$ws.Range("M113").Value = -6827.75
$ws.Range("N113").Value = -13724.3999

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 364.9
$ws.Range("I97").Value = 226.21053
$ws.Range("K97").Value = 226.21053
$ws.Range("M97").Value = 269.78947

# LTW row 11
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 19725
$ws.Range("I11").Value = 19725
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 19725
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -19585
$ws.Range("N11").ClearContents()

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1174.7273
$ws.Range("I81").Value = 731.8570999999999
$ws.Range("K81").Value = 1463.7142
$ws.Range("M81").Value = -402.7141999999999

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1174.7273
$ws.Range("I84").Value = 731.8570999999999
$ws.Range("K84").Value = 7318.571
$ws.Range("M84").Value = -2014.571

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10926.467
$ws.Range("I136").Value = 14508.556
$ws.Range("K136").Value = 43525.66800000001
$ws.Range("M136").Value = -40975.66800000001
